$d = $word.ActiveDocument

$d.Content.Find.Execute("49×38=1862", $true, $true, $false, $false, $false, $true, 1, $false, "75×67=5025", 2) | Out-Null
$d.Content.Find.Execute("72×25=1800", $true, $true, $false, $false, $false, $true, 1, $false, "32×46=1472", 2) | Out-Null
$d.Content.Find.Execute("25×83=2075", $true, $true, $false, $false, $false, $true, 1, $false, "12×87=1044", 2) | Out-Null
$d.Content.Find.Execute("32×27=864", $true, $true, $false, $false, $false, $true, 1, $false, "33×79=2607", 2) | Out-Null
$d.Content.Find.Execute("11×57=627", $true, $true, $false, $false, $false, $true, 1, $false, "27×18=486", 2) | Out-Null
$d.Content.Find.Execute("26×78=2028", $true, $true, $false, $false, $false, $true, 1, $false, "81×45=3645", 2) | Out-Null
$d.Content.Find.Execute("48×48=2304", $true, $true, $false, $false, $false, $true, 1, $false, "27×42=1134", 2) | Out-Null
$d.Content.Find.Execute("18×68=1224", $true, $true, $false, $false, $false, $true, 1, $false, "17×64=1088", 2) | Out-Null
$d.Content.Find.Execute("15×16=240", $true, $true, $false, $false, $false, $true, 1, $false, "35×18=630", 2) | Out-Null
$d.Content.Find.Execute("79×96=7584", $true, $true, $false, $false, $false, $true, 1, $false, "46×77=3542", 2) | Out-Null
$d.Content.Find.Execute("21×27=567", $true, $true, $false, $false, $false, $true, 1, $false, "31×61=1891", 2) | Out-Null
$d.Content.Find.Execute("75×26=1950", $true, $true, $false, $false, $false, $true, 1, $false, "29×33=957", 2) | Out-Null
$d.Content.Find.Execute("30×21=630", $true, $true, $false, $false, $false, $true, 1, $false, "51×59=3009", 2) | Out-Null
$d.Content.Find.Execute("86×34=2924", $true, $true, $false, $false, $false, $true, 1, $false, "73×31=2263", 2) | Out-Null
$d.Content.Find.Execute("88×71=6248", $true, $true, $false, $false, $false, $true, 1, $false, "24×52=1248", 2) | Out-Null
$d.Content.Find.Execute("99×24=2376", $true, $true, $false, $false, $false, $true, 1, $false, "72×62=4464", 2) | Out-Null
$d.Content.Find.Execute("28×73=2044", $true, $true, $false, $false, $false, $true, 1, $false, "99×45=4455", 2) | Out-Null
$d.Content.Find.Execute("47×50=2350", $true, $true, $false, $false, $false, $true, 1, $false, "73×48=3504", 2) | Out-Null
$d.Content.Find.Execute("91×65=5915", $true, $true, $false, $false, $false, $true, 1, $false, "34×65=2210", 2) | Out-Null
$d.Content.Find.Execute("92×14=1288", $true, $true, $false, $false, $false, $true, 1, $false, "82×18=1476", 2) | Out-Null
$d.Content.Find.Execute("18×27=486", $true, $true, $false, $false, $false, $true, 1, $false, "34×77=2618", 2) | Out-Null
$d.Content.Find.Execute("71×25=1775", $true, $true, $false, $false, $false, $true, 1, $false, "50×81=4050", 2) | Out-Null
$d.Content.Find.Execute("69×34=2346", $true, $true, $false, $false, $false, $true, 1, $false, "92×50=4600", 2) | Out-Null
$d.Content.Find.Execute("16×90=1440", $true, $true, $false, $false, $false, $true, 1, $false, "68×32=2176", 2) | Out-Null
$d.Content.Find.Execute("31×95=2945", $true, $true, $false, $false, $false, $true, 1, $false, "24×14=336", 2) | Out-Null
